# Refactor cross_entropy_loss method in SoftmaxBackPropagation class
# Apply targeted cell-value updates to Sheet1 to match the new diagnosis
# output produced by the refactored cross_entropy_loss computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.01

# Row 3
$ws.Range("H3").Value = 0.01

# Row 4
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0

# Row 5
$ws.Range("H5").Value = 0.01

# Row 6
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0.01
